$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General")
$wsLecturer = $wb.Worksheets.Item("Lecturer")
$wsStudents = $wb.Worksheets.Item("Students")

# --- Lecturer sheet: data first (so shared-string order matches) ---
$wsLecturer.Range("A4").Value = "Steven Hawkins"

# --- General sheet ---
$wsGeneral.Range("A1").Value = "Start Date"
$wsGeneral.Range("B1").Value = "End Date"
$wsGeneral.Columns.Item(1).ColumnWidth = 17.5703125
$wsGeneral.Range("A2").Select()

# --- Lecturer sheet header ---
$wsLecturer.Range("A1").Value = "Name"
$wsLecturer.Range("B1").Value = "Datum"

# --- Lecturer sheet rows ---
$wsLecturer.Range("A2").Value = "Albert Einstein"
$wsLecturer.Range("B2").Value = 44806
$wsLecturer.Range("B2").NumberFormat = "mm-dd-yy"
$wsLecturer.Range("C2").Value = 44808
$wsLecturer.Range("C2").NumberFormat = "mm-dd-yy"

$wsLecturer.Range("A3").Value = "Albert Zweistein"

$wsLecturer.Range("B4").Value = 44814
$wsLecturer.Range("B4").NumberFormat = "mm-dd-yy"

$wsLecturer.Range("A5").Value = "Konrad Zuse"
$wsLecturer.Range("B5").Value = 44813
$wsLecturer.Range("B5").NumberFormat = "mm-dd-yy"

$wsLecturer.Range("A6").Value = "Alan Turing"

# Header formatting
$headerRange = $wsLecturer.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.Borders.Item(9).LineStyle = 1
$headerRange.Borders.Item(9).Weight = 2

$wsLecturer.Columns.Item(1).ColumnWidth = 22.140625

# Freeze panes on Lecturer sheet
$wsLecturer.Activate()
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true
$wsLecturer.Range("A8").Select()

# --- Students sheet: clear stray selection (handled automatically by not selecting) ---

$wsLecturer.Activate()
